# Updated cryptos list on Tue Mar 19 04:41:47 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns stay as text (they hold values like "1.00" and
# percentages padded with spaces) instead of being re-interpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value  = "65.662.29"
$ws.Range("E2").Value  = "  -3.78%  "

$ws.Range("D3").Value  = "3.414.30"
$ws.Range("E3").Value  = "  -5.46%  "

$ws.Range("D4").Value  = "1.00"
$ws.Range("E4").Value  = "  +0.25%  "

$ws.Range("D5").Value  = "187.99"
$ws.Range("E5").Value  = "  -7.99%  "

$ws.Range("D6").Value  = "535.34"
$ws.Range("E6").Value  = "  -5.94%  "

$ws.Range("D7").Value  = "0.612"
$ws.Range("E7").Value  = "  -1.48%  "

$ws.Range("D8").Value  = "3.405.14"
$ws.Range("E8").Value  = "  -5.56%  "

$ws.Range("E9").Value  = "  +0.10%  "

$ws.Range("E10").Value = "  -5.78%  "

$ws.Range("D11").Value = "59.91"
$ws.Range("E11").Value = "  -1.22%  "

$ws.Range("E12").Value = "  -11.11%  "

$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  -10.45%  "

$ws.Range("D14").Value = "9.38"
$ws.Range("E14").Value = "  -6.77%  "

$ws.Range("D15").Value = "3.959.98"
$ws.Range("E15").Value = "  -5.18%  "

$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").Value = "3.408.81"
$ws.Range("E17").Value = "  -5.19%  "

$ws.Range("D18").Value = "65.468.55"
$ws.Range("E18").Value = "  -3.66%  "

$ws.Range("D19").Value = "17.69"
$ws.Range("E19").Value = "  -7.18%  "

$ws.Range("D20").Value = "11.31"
$ws.Range("E20").Value = "  -8.74%  "

$ws.Range("D21").Value = "0.987"
$ws.Range("E21").Value = "  -8.26%  "

$ws.Range("D22").Value = "376.91"
$ws.Range("E22").Value = "  -6.55%  "

$ws.Range("D23").Value = "82.66"
$ws.Range("E23").Value = "  -3.15%  "

$ws.Range("D24").Value = "3.80"
$ws.Range("E24").Value = "  -8.96%  "

$ws.Range("D25").Value = "10.99"
$ws.Range("E25").Value = "  -14.37%  "

$ws.Range("D26").Value = "3.71"
$ws.Range("E26").Value = "  -4.84%  "

$ws.Range("D27").Value = "11.82"
$ws.Range("E27").Value = "  -6.06%  "

$ws.Range("D28").Value = "2.69"
$ws.Range("E28").Value = "  -7.85%  "

$ws.Range("D29").Value = "8.64"
$ws.Range("E29").Value = "  -7.95%  "

$ws.Range("D30").Value = "698.15"
$ws.Range("E30").Value = "  +4.30%  "

$ws.Range("D31").Value = "29.95"
$ws.Range("E31").Value = "  -5.34%  "

$ws.Range("E32").Value = "  -17.19%  "

$ws.Range("E33").Value = "  -7.23%  "

$ws.Range("D34").Value = "61.56"
$ws.Range("E34").Value = "  -3.70%  "

$ws.Range("E35").Value = "  -6.11%  "

$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").Value = "37.02"
$ws.Range("E37").Value = "  -12.60%  "

$ws.Range("E38").Value = "  -8.17%  "

$ws.Range("D39").Value = "0.997"
$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("E40").Value = "  -5.05%  "

$ws.Range("D41").Value = "28.53"
$ws.Range("E41").Value = "  +27.57%  "

$ws.Range("D42").Value = "2.908.62"
$ws.Range("E42").Value = "  -10.90%  "

$ws.Range("E43").Value = "  -10.95%  "

$ws.Range("D44").Value = "0.0405"
$ws.Range("E44").Value = "  -3.19%  "

$ws.Range("E45").Value = "  -3.74%  "

$ws.Range("D46").Value = "0.0₃0631"
$ws.Range("E46").Value = "  -17.80%  "

$ws.Range("D47").Value = "2.39"
$ws.Range("E47").Value = "  -12.86%  "

$ws.Range("E48").Value = "  -3.31%  "

$ws.Range("D49").Value = "137.92"
$ws.Range("E49").Value = "  -0.77%  "

# Row 50 and 51 swap their coin (Stacks <-> ApeXProtocol) along with new figures.
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "2.92"
$ws.Range("E50").Value = "  -5.00%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "2.67"
$ws.Range("E51").Value = "  -2.82%  "
